# (hel-884) Correction du wording sur le taux de vétusté de constructions
#
# Fix the typo "construction" -> "constructions" in the label of the
# "Taux de vétusté des construction(s) (en %)" indicator, on the
# "Lisez-moi" sheet, cell A18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lisez-moi")

$ws.Range("A18").Value = "Taux de vétusté des constructions (en %)"
